$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.368.48'
$ws.Range("E2").Value = '  +0.52%  '

$ws.Range("D3").Value = '2.521.37'
$ws.Range("E3").Value = '  +2.90%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '543.86'
$ws.Range("E5").Value = '  +0.69%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.00'
$ws.Range("E6").Value = '  -1.60%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.574'
$ws.Range("E8").Value = '  +0.51%  '

$ws.Range("D9").Value = '2.553.47'
$ws.Range("E9").Value = '  +3.54%  '

$ws.Range("E10").Value = '  +2.08%  '

$ws.Range("E11").Value = '  +0.58%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.62'
$ws.Range("E12").Value = '  +5.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.356'
$ws.Range("E13").Value = '  +1.19%  '

$ws.Range("D14").Value = '2.999.87'
$ws.Range("E14").Value = '  +3.88%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.87'
$ws.Range("E15").Value = '  -0.69%  '

$ws.Range("D16").Value = '59.342.86'
$ws.Range("E16").Value = '  +0.69%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000141'
$ws.Range("E17").Value = '  +2.54%  '

$ws.Range("D18").Value = '2.542.34'
$ws.Range("E18").Value = '  +1.06%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.28'
$ws.Range("E19").Value = '  +1.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.31'
$ws.Range("E20").Value = '  -0.96%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '327.08'
$ws.Range("E21").Value = '  +0.66%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  +3.46%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.87'
$ws.Range("E23").Value = '  +2.76%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.09'
$ws.Range("E24").Value = '  +2.27%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.440'
$ws.Range("E25").Value = '  -4.16%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.164'
$ws.Range("E26").Value = '  +2.47%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.993'
$ws.Range("E27").Value = '  +1.60%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.00'
$ws.Range("E28").Value = '  +3.80%  '

$ws.Range("B29").Value = 'Aptos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.83'
$ws.Range("E29").Value = '  +1.95%  '

$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0786'
$ws.Range("E30").Value = '  +1.67%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.83'
$ws.Range("E31").Value = '  +0.36%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.23'
$ws.Range("E32").Value = '  -3.62%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.50'
$ws.Range("E33").Value = '  +9.44%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").Value = '  +0.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '156.64'
$ws.Range("E35").Value = '  +0.02%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.75'
$ws.Range("E36").Value = '  +1.70%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.41'
$ws.Range("E37").Value = '  -1.27%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.63'
$ws.Range("E38").Value = '  -4.89%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.70'
$ws.Range("E39").Value = '  -2.70%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.03'
$ws.Range("E40").Value = '  +2.27%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '300.90'
$ws.Range("E41").Value = '  -5.29%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.73'
$ws.Range("E42").Value = '  +0.36%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.828'
$ws.Range("E43").Value = '  -1.18%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.992'
$ws.Range("E44").Value = '  -0.35%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.606'
$ws.Range("E45").Value = '  +4.43%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.80'
$ws.Range("E46").Value = '  +0.60%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0934'
$ws.Range("E47").Value = '  -0.95%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.82'
$ws.Range("E48").Value = '  +2.36%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.02'
$ws.Range("E49").Value = '  +1.95%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0229'
$ws.Range("E50").Value = '  -0.42%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0516'
$ws.Range("E51").Value = '  -1.84%  '
